# UC001 - Calcular Desconto de Produto
# "Alterada a etapa 12" -> on closer inspection this commit actually updates
# step 7 (row 92, TC5) plus the shared "Tipo de Cliente" expected-result text
# that is reused by step 6 of every test case (TC1..TC5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 6 "Expected Results" text is shared by TC1 (D15), TC2 (D34), TC3 (D53),
# TC4 (D72) and TC5 (D91): drop the "(padrão 'A')" / "(vazio)" qualifiers.
$newExpected = "SYSTEM apresenta campos: Tipo de Cliente e Quantidade"
$ws.Range("D15").Value = $newExpected
$ws.Range("D34").Value = $newExpected
$ws.Range("D53").Value = $newExpected
$ws.Range("D72").Value = $newExpected
$ws.Range("D91").Value = $newExpected

# TC5 step 7 (row 92): the step now changes the client type to B (not C).
$ws.Range("B92").Value = "Usuário do Sistema altera para tipo de cliente B"
